$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.925.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.23%  "

# Row 3
$ws.Range("D3").Value = "'3.268.31"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'586.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.34%  "

# Row 6
$ws.Range("D6").Value = "'186.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.83%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.600"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.38%  "

# Row 9
$ws.Range("E9").Value = "  +3.47%  "

# Row 10
$ws.Range("D10").Value = "'6.71"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.78%  "

# Row 11
$ws.Range("E11").Value = "  +0.65%  "

# Row 12
$ws.Range("D12").Value = "'3.840.55"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.67%  "

# Row 13
$ws.Range("E13").Value = "  +0.56%  "

# Row 14
$ws.Range("D14").Value = "'28.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.64%  "

# Row 15
$ws.Range("D15").Value = "'67.962.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.29%  "

# Row 16
$ws.Range("E16").Value = "  +2.27%  "

# Row 17
$ws.Range("D17").Value = "'3.275.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.61%  "

# Row 18
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").Value = "'13.64"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.64%  "

# Row 20
$ws.Range("D20").Value = "'382.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.61%  "

# Row 21
$ws.Range("D21").Value = "'7.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.95%  "

# Row 22
$ws.Range("D22").Value = "'71.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("D24").Value = "'0.515"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$ws.Range("E25").Value = "  +1.49%  "

# Row 26
$ws.Range("D26").Value = "'0.189"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.75%  "

# Row 27
$ws.Range("E27").Value = "  -0.54%  "

# Row 28
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("E29").Value = "  +4.14%  "

# Row 30
$ws.Range("E30").Value = "  +0.92%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'22.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.12%  "

# Row 32
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'7.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.83%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  +0.53%  "

# Row 35
$ws.Range("E35").Value = "  +2.40%  "

# Row 36
$ws.Range("D36").Value = "'162.24"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.95%  "

# Row 37
$ws.Range("D37").Value = "'1.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.00%  "

# Row 38
$ws.Range("D38").Value = "'0.838"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.06%  "

# Row 39
$ws.Range("D39").Value = "'6.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.42%  "

# Row 40
$ws.Range("D40").Value = "'26.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.30%  "

# Row 41
$ws.Range("E41").Value = "  +4.64%  "

# Row 42
$ws.Range("E42").Value = "  +0.56%  "

# Row 43
$ws.Range("D43").Value = "'0.0693"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.56%  "

# Row 44
$ws.Range("D44").Value = "'41.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.03%  "

# Row 45
$ws.Range("D45").Value = "'25.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").Value = "'2.648.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.05%  "

# Row 47
$ws.Range("D47").Value = "'341.97"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.64%  "

# Row 48
$ws.Range("D48").Value = "'0.0284"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.66%  "

# Row 49
$ws.Range("D49").Value = "'31.98"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.18%  "

# Row 50
$ws.Range("D50").Value = "'0.997"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.11%  "

# Row 51
$ws.Range("E51").Value = "  -0.33%  "
